$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Requisitos" block occupies rows 23-25 in columns B and C.
# Originally:
#   Row 23: LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)
#   Row 24: LOB1021 -  Física IV  (Requisito)
#   Row 25: LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)
# The diff moves the LOM3246 entry to the end, shifting the other two up:
#   Row 23: LOB1021 -  Física IV  (Requisito)
#   Row 24: LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)
#   Row 25: LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)

$lom3246 = "LOM3246 -  Técnicas de Caracterização de Materiais  (Indicação de Conjunto)`n"
$lob1021 = "LOB1021 -  Física IV  (Requisito)`n"
$lom3016 = "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)`n"

$ws.Range("B23").Value = $lob1021
$ws.Range("C23").Value = $lob1021

$ws.Range("B24").Value = $lom3016
$ws.Range("C24").Value = $lom3016

$ws.Range("B25").Value = $lom3246
$ws.Range("C25").Value = $lom3246
